# Refresh Leve profit figures (currentAveragePrice*, LevePrice*, LeveProfit*)
# across the per-job sheets, as produced by the scheduled market-data runner.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 3
$ws.Range("H3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("L3").ClearContents()
$ws.Range("N3").Value = 0

# Row 40
$ws.Range("H40").Value = 2253.6155
$ws.Range("I40").Value = 1964.1428
$ws.Range("J40").Value = 2591.3333
$ws.Range("K40").Value = 1964.1428
$ws.Range("L40").Value = 2591.3333
$ws.Range("M40").Value = -1789.1428
$ws.Range("N40").Value = -2941.3333

# Row 74
$ws.Range("H74").Value = 8534.333000000001
$ws.Range("I74").Value = 2801.5
$ws.Range("K74").Value = 2801.5
$ws.Range("M74").Value = -1865.5

# Row 76
$ws.Range("H76").Value = 6358.5
$ws.Range("I76").Value = 5464.6665
$ws.Range("K76").Value = 5464.6665
$ws.Range("M76").Value = -5149.6665

# Row 77
$ws.Range("H77").Value = 8534.333000000001
$ws.Range("I77").Value = 2801.5
$ws.Range("K77").Value = 14007.5
$ws.Range("M77").Value = -9327.5

# Row 79
$ws.Range("H79").Value = 6358.5
$ws.Range("I79").Value = 5464.6665
$ws.Range("K79").Value = 5464.6665
$ws.Range("M79").Value = -4372.6665

# Row 88
$ws.Range("H88").Value = 2108.2222
$ws.Range("J88").Value = 2396.8333
$ws.Range("L88").Value = 2396.8333
$ws.Range("N88").Value = -3208.8333

# Row 91
$ws.Range("H91").Value = 2108.2222
$ws.Range("J91").Value = 2396.8333
$ws.Range("L91").Value = 2396.8333
$ws.Range("N91").Value = -5204.8333

# Row 102
$ws.Range("H102").Value = 0
$ws.Range("J102").Value = 0
$ws.Range("L102").ClearContents()
$ws.Range("N102").Value = 0

# Row 131
$ws.Range("H131").Value = 5385.3076
$ws.Range("I131").Value = 1429.8572
$ws.Range("K131").Value = 4289.571599999999
$ws.Range("M131").Value = 750.4284000000007

# Row 141
$ws.Range("H141").Value = 2273.75
$ws.Range("I141").Value = 2273.75
$ws.Range("K141").Value = 6821.25
$ws.Range("M141").Value = -1641.25

$ws = $wb.Worksheets.Item("ARM")
# Row 88
$ws.Range("H88").Value = 4001.6667
$ws.Range("I88").Value = 4002.5
$ws.Range("J88").Value = 4000
$ws.Range("K88").Value = 4002.5
$ws.Range("L88").Value = 4000
$ws.Range("M88").Value = -3596.5
$ws.Range("N88").Value = -4812

# Row 91
$ws.Range("H91").Value = 4001.6667
$ws.Range("I91").Value = 4002.5
$ws.Range("J91").Value = 4000
$ws.Range("K91").Value = 4002.5
$ws.Range("L91").Value = 4000
$ws.Range("M91").Value = -2598.5
$ws.Range("N91").Value = -6808

$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 1764.8334
$ws.Range("I20").Value = 1724.25
$ws.Range("K20").Value = 1724.25
$ws.Range("M20").Value = -1477.25

# Row 82
$ws.Range("H82").Value = 29948.5
$ws.Range("I82").Value = 6598.3335
$ws.Range("J82").Value = 99999
$ws.Range("K82").Value = 6598.3335
$ws.Range("L82").Value = 99999
$ws.Range("M82").Value = -6215.3335
$ws.Range("N82").Value = -100765

# Row 85
$ws.Range("H85").Value = 29948.5
$ws.Range("I85").Value = 6598.3335
$ws.Range("J85").Value = 99999
$ws.Range("K85").Value = 6598.3335
$ws.Range("L85").Value = 99999
$ws.Range("M85").Value = -5272.3335
$ws.Range("N85").Value = -102651

# Row 86
$ws.Range("H86").Value = 3300.3333
$ws.Range("I86").Value = 2500
$ws.Range("K86").Value = 2500
$ws.Range("M86").Value = -1377

# Row 89
$ws.Range("H89").Value = 3300.3333
$ws.Range("I89").Value = 2500
$ws.Range("K89").Value = 12500
$ws.Range("M89").Value = -6884

# Row 105
$ws.Range("H105").Value = 4361.926
$ws.Range("I105").Value = 3642.0625
$ws.Range("K105").Value = 3642.0625
$ws.Range("M105").Value = -1895.0625

$ws = $wb.Worksheets.Item("CRP")
# Row 62
$ws.Range("H62").Value = 32302.285
$ws.Range("J62").Value = 70414.836
$ws.Range("L62").Value = 70414.836
$ws.Range("N62").Value = -71662.836

# Row 65
$ws.Range("H65").Value = 32302.285
$ws.Range("J65").Value = 70414.836
$ws.Range("L65").Value = 352074.18
$ws.Range("N65").Value = -358314.18

# Row 134
$ws.Range("H134").Value = 2535.65
$ws.Range("I134").Value = 1407.1428
$ws.Range("J134").Value = 5168.8335
$ws.Range("K134").Value = 4221.428400000001
$ws.Range("L134").Value = 15506.5005
$ws.Range("M134").Value = -1686.428400000001
$ws.Range("N134").Value = -20576.5005

$ws = $wb.Worksheets.Item("CUL")
# Row 49
$ws.Range("H49").Value = 925.5
$ws.Range("I49").Value = 399
$ws.Range("J49").Value = 1452
$ws.Range("K49").Value = 1197
$ws.Range("L49").Value = 4356
$ws.Range("M49").Value = -1041
$ws.Range("N49").Value = -4668

# Row 113
$ws.Range("H113").Value = 972.3889
$ws.Range("I113").Value = 1757.6
$ws.Range("J113").Value = 845.74194
$ws.Range("K113").Value = 5272.799999999999
$ws.Range("L113").Value = 2537.22582
$ws.Range("M113").Value = -3102.799999999999
$ws.Range("N113").Value = -6877.22582

# Row 127
$ws.Range("H127").Value = 2000
$ws.Range("J127").Value = 2000
$ws.Range("L127").Value = 6000
$ws.Range("N127").Value = -15920

$ws = $wb.Worksheets.Item("GSM")
# Row 20
$ws.Range("H20").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("L20").ClearContents()
$ws.Range("N20").Value = 0

# Row 70
$ws.Range("H70").Value = 6891.3335
$ws.Range("I70").Value = 6100.857
$ws.Range("K70").Value = 6100.857
$ws.Range("M70").Value = -5830.857

# Row 73
$ws.Range("H73").Value = 6891.3335
$ws.Range("I73").Value = 6100.857
$ws.Range("K73").Value = 6100.857
$ws.Range("M73").Value = -5164.857

# Row 80
$ws.Range("H80").Value = 6031
$ws.Range("I80").Value = 4487.857
$ws.Range("K80").Value = 4487.857
$ws.Range("M80").Value = -3489.857

# Row 83
$ws.Range("H83").Value = 6031
$ws.Range("I83").Value = 4487.857
$ws.Range("K83").Value = 22439.285
$ws.Range("M83").Value = -17447.285

# Row 122
$ws.Range("H122").Value = 649858.25
$ws.Range("I122").Value = 85383.164
$ws.Range("J122").Value = 2004598.4
$ws.Range("K122").Value = 256149.492
$ws.Range("L122").Value = 6013795.199999999
$ws.Range("M122").Value = -253699.492
$ws.Range("N122").Value = -6018695.199999999

# Row 136
$ws.Range("H136").Value = 26854.215
$ws.Range("J136").Value = 26854.215
$ws.Range("L136").Value = 80562.645
$ws.Range("N136").Value = -85662.645

$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 647.6667
$ws.Range("I22").Value = 845.8
$ws.Range("J22").Value = 400
$ws.Range("K22").Value = 845.8
$ws.Range("L22").Value = 400
$ws.Range("M22").Value = -550.8
$ws.Range("N22").Value = -990

# Row 27
$ws.Range("H27").Value = 647.6667
$ws.Range("I27").Value = 845.8
$ws.Range("J27").Value = 400
$ws.Range("K27").Value = 845.8
$ws.Range("L27").Value = 400
$ws.Range("M27").Value = -738.8
$ws.Range("N27").Value = -614

# Row 46
$ws.Range("H46").Value = 2959.0908
$ws.Range("I46").Value = 2378.5715
$ws.Range("K46").Value = 2378.5715
$ws.Range("M46").Value = -2190.5715

# Row 55
$ws.Range("H55").Value = 417.66666
$ws.Range("I55").Value = 377
$ws.Range("J55").Value = 560
$ws.Range("K55").Value = 377
$ws.Range("L55").Value = 560
$ws.Range("M55").Value = -204
$ws.Range("N55").Value = -906

# Row 61
$ws.Range("H61").Value = 3862.2144
$ws.Range("I61").Value = 3839.0833
$ws.Range("K61").Value = 3839.0833
$ws.Range("M61").Value = -3637.0833

# Row 113
$ws.Range("H113").Value = 3862.2144
$ws.Range("I113").Value = 3839.0833
$ws.Range("K113").Value = 3839.0833
$ws.Range("M113").Value = -1669.0833

# Row 118
$ws.Range("H118").Value = 31000
$ws.Range("J118").Value = 31000
$ws.Range("L118").Value = 31000
$ws.Range("N118").Value = -34314

# Row 122
$ws.Range("H122").Value = 2998.5
$ws.Range("I122").Value = 2998.5
$ws.Range("K122").Value = 8995.5
$ws.Range("M122").Value = -6545.5

$ws = $wb.Worksheets.Item("WVR")
# Row 7
$ws.Range("H7").Value = 1089.8334
$ws.Range("I7").Value = 1478
$ws.Range("J7").Value = 701.6667
$ws.Range("K7").Value = 1478
$ws.Range("L7").Value = 701.6667
$ws.Range("M7").Value = -1365
$ws.Range("N7").Value = -927.6667

# Row 15
$ws.Range("H15").Value = 20997.5
$ws.Range("I15").Value = 20995
$ws.Range("J15").Value = 21000
$ws.Range("K15").Value = 20995
$ws.Range("L15").Value = 21000
$ws.Range("M15").Value = -20707
$ws.Range("N15").Value = -21576

# Row 20
$ws.Range("H20").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("L20").ClearContents()
$ws.Range("N20").Value = 0

# Row 122
$ws.Range("H122").Value = 1949.5
$ws.Range("I122").Value = 1949.5
$ws.Range("K122").Value = 5848.5
$ws.Range("M122").Value = -3398.5

# Row 136
$ws.Range("H136").Value = 30863.117
$ws.Range("I136").Value = 1233.9546
$ws.Range("J136").Value = 85183.25
$ws.Range("K136").Value = 3701.8638
$ws.Range("L136").Value = 255549.75
$ws.Range("M136").Value = -1151.8638
$ws.Range("N136").Value = -260649.75
